$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("Estimated time" etc. shift right by one)
$ws.Columns("B:B").Insert()

# Header for the new Priority column
$ws.Range("B1").Value = "Priority"

# Seed the new shared strings in the same first-use order as the authored
# workbook (Priority, Med, High, Low, Will, Liam) before filling in the rest.
$ws.Range("B7").Value = "Med"
$ws.Range("B2").Value = "High"
$ws.Range("B4").Value = "Low"
$ws.Range("D2").Value = "Will"
$ws.Range("D3").Value = "Liam"

# Priority column (B)
$ws.Range("B3").Value = "High"
$ws.Range("B5").Value = "High"
$ws.Range("B6").Value = "High"
$ws.Range("B8").Value = "Med"
$ws.Range("B9").Value = "Med"
$ws.Range("B10").Value = "Med"
$ws.Range("B11").Value = "Med"
$ws.Range("B12").Value = "Low"

# Estimate column (C), in sprints
$ws.Range("C2").Value = 3
$ws.Range("C3").Value = 3
$ws.Range("C4").Value = 3
$ws.Range("C5").Value = 4
$ws.Range("C6").Value = 4
$ws.Range("C7").Value = 4
$ws.Range("C8").Value = 5
$ws.Range("C9").Value = 5
$ws.Range("C10").Value = 3
$ws.Range("C11").Value = 4
$ws.Range("C12").Value = 5

# Person assigned (D) -- only filled in for sprint 1 stories
$ws.Range("D4").Value = "Liam"
$ws.Range("D5").Value = "Will"

# Match the narrower width used for the new Priority column
$ws.Columns("B:B").ColumnWidth = 6

# Selection / active cell as recorded at save time
[void]$ws.Range("D6").Select()
